# Apply the "Updated cryptos list" refresh: new Price (D) and Volume(1h) (E)
# values for the coinranking snapshot rows (row 2..51 of Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.015.52"
$ws.Range("E2").Value = "  +3.84%  "
$ws.Range("D3").Value = "3.636.52"
$ws.Range("E3").Value = "  +6.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.52"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.51"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "3.628.25"
$ws.Range("E7").Value = "  +6.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.612"
$ws.Range("E8").Value = "  +1.95%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.204"
$ws.Range("E10").Value = "  +3.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.608"
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "50.08"
$ws.Range("E12").Value = "  +2.98%  "
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "693.71"
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").Value = "4.227.52"
$ws.Range("E15").Value = "  +6.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "9.02"
$ws.Range("E16").Value = "  +3.83%  "
$ws.Range("D17").Value = "3.689.47"
$ws.Range("E17").Value = "  +8.13%  "
$ws.Range("D18").Value = "71.965.54"
$ws.Range("E18").Value = "  +3.57%  "
$ws.Range("E19").Value = "  +2.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.41"
$ws.Range("E20").Value = "  +3.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.62"
$ws.Range("E21").Value = "  +2.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.938"
$ws.Range("E22").Value = "  +2.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.81"
$ws.Range("E23").Value = "  +7.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.95"
$ws.Range("E24").Value = "  +4.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "104.16"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.04"
$ws.Range("E26").Value = "  +2.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.88"
$ws.Range("E27").Value = "  +5.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("E28").Value = "  +3.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.22"
$ws.Range("E29").Value = "  +3.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.16"
$ws.Range("E30").Value = "  +3.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.34"
$ws.Range("E31").Value = "  +5.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.17"
$ws.Range("E32").Value = "  +15.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "584.90"
$ws.Range("E33").Value = "  +4.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.39"
$ws.Range("E34").Value = "  +1.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.110"
$ws.Range("E35").Value = "  +2.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.69"
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "3.664.49"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.03"
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("D41").Value = "0.0₃0770"
$ws.Range("E41").Value = "  +6.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.44"
$ws.Range("E42").Value = "  +5.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0463"
$ws.Range("E43").Value = "  +8.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.78"
$ws.Range("E44").Value = "  +3.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.351"
$ws.Range("E45").Value = "  +3.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.43"
$ws.Range("E46").Value = "  +1.94%  "
$ws.Range("E47").Value = "  +6.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.134"
$ws.Range("E48").Value = "  +2.83%  "
$ws.Range("E49").Value = "  +3.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.35"
$ws.Range("E51").Value = "  -0.28%  "
